$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-11: column B becomes 1, column C becomes 0.7 ---
# (A stays as the running index 1..10, D2 formula/E2 value also change)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
    $ws.Cells.Item($r, 3).Value = 0.7
}

# E2 (k) changes from 20 to 5
$ws.Range("E2").Value = 5

# --- Append new rows 12-21 with the same pattern: index, weight=1, reliability=0.7 ---
for ($r = 12; $r -le 21; $r++) {
    $idx = $r - 1
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = 1
    $ws.Cells.Item($r, 3).Value = 0.7
}

# --- Column C (reliability) width (closest representable value to 10.81640625) ---
$ws.Columns.Item(3).ColumnWidth = 10.81640625

# --- Turn A1:E33 into a table named "Tabell1" ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:E33"), $null, 1)
$tbl.Name = "Tabell1"

# --- Selection moves to B6 ---
$ws.Range("B6").Select()
